$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Strip the _GoBack bookmark from the "Next week: Finish
#    administrator Panel..." paragraph - it will be re-added further
#    down, attached to the new "We have fixed errors." paragraph.
# ------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Find the anchor paragraph by its text (Paragraph.Range.Text includes
# the trailing paragraph mark, so use Contains rather than equality).
$anchorText = "Next week: Finish administrator Panel, Keep working on JPA manager and Manager."
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Contains($anchorText)) {
        $anchorIndex = $i
    }
}

# The paragraph right after the anchor is an already-empty paragraph
# that must stay untouched; the second empty paragraph after that is
# the one that receives the "We have fixed errors." text.
$emptyIndex1 = $anchorIndex + 1
$emptyIndex2 = $anchorIndex + 2

# ------------------------------------------------------------------
# 2. Create all the new paragraph marks FIRST, before applying any
#    character formatting - InsertParagraphAfter() copies the
#    formatting of the mark it splits off, so doing this up-front
#    keeps the underline used for the "Friday, 20/04/2018" heading
#    from leaking into the paragraphs that follow it.
# ------------------------------------------------------------------
$pEmpty1 = $d.Paragraphs($emptyIndex1)
$pEmpty1.Range.InsertParagraphAfter()
$fridayIndex = $emptyIndex1 + 1

$pFriday = $d.Paragraphs($fridayIndex)
$pFriday.Range.InsertParagraphAfter()
$thisWeekIndex = $fridayIndex + 1

# $emptyIndex2 has shifted by the two paragraphs inserted above.
$errorsIndex = $emptyIndex2 + 2
$pErrors = $d.Paragraphs($errorsIndex)
$pErrors.Range.InsertParagraphAfter()
$nextWeekIndex = $errorsIndex + 1

# ------------------------------------------------------------------
# 3. "Friday, 20/04/2018" heading (underlined).
# ------------------------------------------------------------------
$pFriday = $d.Paragraphs($fridayIndex)
$pFriday.Range.Text = "Friday, 20/04/2018"
$pFriday.Range.Font.Underline = 1

# ------------------------------------------------------------------
# 4. "This week: ..." paragraph (three runs with identical
#    formatting, built up as separate InsertAfter calls).
# ------------------------------------------------------------------
$pThisWeek = $d.Paragraphs($thisWeekIndex)
$pThisWeek.Range.Text = "This week: "
$r = $d.Range($pThisWeek.Range.End - 1, $pThisWeek.Range.End - 1)
$r.InsertAfter("We did the setters and getters of the Packaged and Arrives POJOs, fixed the JPA annotations and the inserts for JPA, we also did many Selects in JPA and the deletes in JPA")
$r = $d.Range($pThisWeek.Range.End - 1, $pThisWeek.Range.End - 1)
$r.InsertAfter(". ")

# ------------------------------------------------------------------
# 5. "We have fixed errors." paragraph, followed by the restored
#    _GoBack bookmark (collapsed) and a trailing space.
#
#    To work around a bug where adding a zero-length bookmark near
#    the tail of the document snaps to the wrong range, a one-char
#    placeholder is bookmarked and then deleted through the bookmark
#    itself, which correctly collapses the bookmark in place.
# ------------------------------------------------------------------
$pErrors = $d.Paragraphs($errorsIndex)
$pErrors.Range.Text = "We have fixed errors.#"
$pErrors.Range.Font.LanguageID = "en-US"
$placeholderPos = $pErrors.Range.End - 2
$bmRange = $d.Range($placeholderPos, $placeholderPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$bm = $d.Bookmarks("_GoBack")
$bm.Range.Text = ""

$r = $d.Range($pErrors.Range.End - 1, $pErrors.Range.End - 1)
$r.InsertAfter(" ")

# ------------------------------------------------------------------
# 6. Final "Next week: ..." paragraph.
# ------------------------------------------------------------------
$pNextWeek = $d.Paragraphs($nextWeekIndex)
$pNextWeek.Range.Text = "Next week:"
$r = $d.Range($pNextWeek.Range.End - 1, $pNextWeek.Range.End - 1)
$r.InsertAfter(" Keep working on JPA manager and SQL manager")
